# "wrapping up test file audit"
#
# The optimization_parameters sheet had a stray leftover test row (row 16,
# labeled "Sheet" with values 3 and 4) that doesn't belong with the real
# optimization parameters. Remove it -- this shifts the
# simulation_timepoints row up to row 16, and Excel cleans up the now-
# unused "Sheet" shared string / number-format style as part of the save.
#
# Finish by leaving the workbook focused back on the first sheet
# (production_rates), matching the window state the file was committed with.

$wb = $excel.ActiveWorkbook

$paramsSheet = $wb.Worksheets.Item("optimization_parameters")
$paramsSheet.Rows(16).Delete() | Out-Null

# Mirror Excel's behavior of moving the selection down onto the row that
# slides up to fill the gap (now the whole row 16).
$paramsSheet.Range("A16:XFD16").Select() | Out-Null

# Return focus to the first sheet before saving, like the committed file.
$wb.Worksheets.Item("production_rates").Select() | Out-Null
